# Updates the "想去人数" (F column) counts on the "展览" and "全部类型"
# worksheets to reflect the latest generated data (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F, applied identically to
# both the "展览" and "全部类型" sheets (they carry duplicate data).
$updates = @{
    2  = 1393
    3  = 2676
    4  = 556
    6  = 6620
    7  = 499
    9  = 12
    11 = 6
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
